$d = $word.ActiveDocument

# The document currently contains three separate runs forming the text:
#   <id>  p126r_a1  </id>
# They need to be merged into a single run reading: <id>p126r_1</id>
# using the formatting of the first run (Courier New, color 7f6000, sz 18).

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("<id>p126r_a1</id>", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "<id>p126r_1</id>", 2)
